$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: 20240304 ---
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "20240304"
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)  # xlPasteFormats - restore General format + style from A6

$ws.Range("B7").Value = 72
$ws.Range("C7").Value = 156
$ws.Range("D7").Value = 69
$ws.Range("E7").Value = 107
$ws.Range("F7").Value = 162
$ws.Range("G7").Value = 63

# --- Row 8: 20240305 ---
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = "20240305"
$ws.Range("A6").Copy()
$ws.Range("A8").PasteSpecial(-4122)  # xlPasteFormats - restore General format + style from A6

$ws.Range("B8").Value = 71
$ws.Range("C8").Value = 339
$ws.Range("D8").Value = 68
$ws.Range("E8").Value = 107
$ws.Range("F8").Value = 167
$ws.Range("G8").Value = 58
